$wb = $excel.ActiveWorkbook

# --- Sheet "Vendas": remove rows 4-7, update remaining data rows ---
$wsVendas = $wb.Worksheets.Item("Vendas")

$wsVendas.Rows("4:7").Delete()

$wsVendas.Range("A2").Value = "Cliente Porta"
$wsVendas.Range("B2").Value = "Cimento Nassau 50kg"
$wsVendas.Range("C2").Value = 35.0
$wsVendas.Range("D2").Value = 1
$wsVendas.Range("E2").Value = "A Vista - Pix"
$wsVendas.Range("F2").Value = "27/04/2023"

$wsVendas.Range("A3").Value = "Cliente Porta"
$wsVendas.Range("B3").Value = "Cimento Nassau 50kg"
$wsVendas.Range("C3").Value = 70.0
$wsVendas.Range("D3").Value = 2
$wsVendas.Range("E3").Value = "A Vista - Pix"
$wsVendas.Range("F3").Value = "26/04/2023"

# --- Sheet "Ganhos": update the summary row ---
$wsGanhos = $wb.Worksheets.Item("Ganhos")

$wsGanhos.Range("A2").Value = 105.0
$wsGanhos.Range("B2").Value = 105.0
$wsGanhos.Range("C2").Value = 105.0

# "04" must stay text (leading zero) - format as Text before assigning,
# otherwise Excel parses it as the number 4.
$wsGanhos.Range("D2").NumberFormat = "@"
$wsGanhos.Range("D2").Value = "04"
